$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124, shifting rows 124-199 down to 125-200.
$ws.Rows(124).Insert()

# Populate the new row 124 with its values.
$ws.Cells.Item(124, 1).Value = 7
$ws.Cells.Item(124, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(124, 3).Value = "Ñuble"
$ws.Cells.Item(124, 4).Value = 44897
$ws.Cells.Item(124, 5).Value = 16
$ws.Cells.Item(124, 6).Value = 100112028
$ws.Cells.Item(124, 7).Value = "Sandia"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 400
$ws.Cells.Item(124, 11).Value = 750
$ws.Cells.Item(124, 12).Value = 800
$ws.Cells.Item(124, 13).Value = 775
$ws.Cells.Item(124, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(124, 15).Value = "Paine"
$ws.Cells.Item(124, 16).Value = 775
$ws.Cells.Item(124, 17).Value = 1
$ws.Cells.Item(124, 18).Value = "Hortaliza"
